$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-name" and "codeforiati:group-code" columns (C and D)
# were reordered so that group-code now comes before group-name: swap the
# contents of columns C and D for the header row and every data row.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $groupNameVal = $ws.Cells.Item($r, 3).Value2
    $groupCodeVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $groupCodeVal
    $ws.Cells.Item($r, 4).Value = $groupNameVal
}
